$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Waves")
$ws.Range("B2").Value = 2.4
